$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original layout: A1=1, A2=2, A3=3, A4=4
# Target layout:   A1="1:4" (text, Courier New), A2=blank, A3=1, A4=2, A5=3, A6=4

# Shift the numbers down two rows (write from bottom to top to avoid overwrite issues).
$ws.Range("A6").Value = $ws.Range("A4").Value2
$ws.Range("A5").Value = $ws.Range("A3").Value2
$ws.Range("A4").Value = $ws.Range("A2").Value2
$ws.Range("A3").Value = $ws.Range("A1").Value2

# Clear A2 (now empty, between label row and numbers).
$ws.Range("A2").ClearContents()

# Set A1 to the text label styled with Courier New font.
$ws.Range("A1").Value = "1:4"
$ws.Range("A1").Font.Name = "Courier New"
